# Cotações atualizadas - 2025-09-07
# Append a new row (row 3) with the quotes for 2025-09-07, reusing the
# same number format as the previous date row so the new cell matches
# the existing "Data" column formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("A3").Value = 45907
$ws.Range("B3").Value = "20,8447"
$ws.Range("C3").Value = "14,7888"
$ws.Range("D3").Value = "14,7856"
$ws.Range("E3").Value = "13,1918"
